$d = $word.ActiveDocument

# Locate the target empty paragraph to replace. Structurally, the document
# body ends with a trailing empty "Sin espaciado" paragraph (immediately
# before the sectPr). The paragraph to replace is the empty paragraph right
# before that trailing one (it carries no visible text, just <w:lang
# w:val="es-ES"/> in its rPr in the original file).
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($count)
$targetIndex = $count - 1
$target = $d.Paragraphs($targetIndex)

# Sanity checks so we fail loudly instead of silently editing the wrong spot.
if ($lastPara.Style.NameLocal -ne "No Spacing") {
    throw "Unexpected document tail: last paragraph style is '$($lastPara.Style.NameLocal)', expected 'No Spacing'"
}
if ($target.Range.Text.Trim() -ne "") {
    throw "Unexpected document tail: target paragraph is not empty"
}

$r = $target.Range

$xml = @'
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="Sinespaciado"/>
        <w:tabs>
          <w:tab w:val="left" w:pos="5670"/>
        </w:tabs>
        <w:spacing w:line="276" w:lineRule="auto"/>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:b/>
          <w:color w:val="4472C4" w:themeColor="accent1"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:b/>
          <w:bCs/>
          <w:noProof/>
          <w:color w:val="4472C4" w:themeColor="accent1"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve">HUJR-001 </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:b/>
          <w:color w:val="4472C4" w:themeColor="accent1"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>Registro de proveedor</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Sinespaciado"/>
        <w:tabs>
          <w:tab w:val="left" w:pos="5670"/>
        </w:tabs>
        <w:spacing w:line="276" w:lineRule="auto"/>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rStyle w:val="Textoennegrita"/>
          <w:rFonts w:cstheme="minorHAnsi"/>
        </w:rPr>
        <w:t>Como</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
        </w:rPr>
        <w:t>: administrador del sistema,</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
        </w:rPr>
        <w:br/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rStyle w:val="Textoennegrita"/>
          <w:rFonts w:cstheme="minorHAnsi"/>
        </w:rPr>
        <w:t>quiero</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rStyle w:val="Textoennegrita"/>
          <w:rFonts w:cstheme="minorHAnsi"/>
        </w:rPr>
        <w:t>:</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
        </w:rPr>
        <w:t xml:space="preserve"> registrar nuevos proveedores con su nombre</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
        </w:rPr>
        <w:t>, teléfono, correo y dirección,</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Sinespaciado"/>
        <w:tabs>
          <w:tab w:val="left" w:pos="5670"/>
        </w:tabs>
        <w:spacing w:line="276" w:lineRule="auto"/>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:b/>
          <w:bCs/>
          <w:noProof/>
          <w:color w:val="4472C4" w:themeColor="accent1"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rStyle w:val="Textoennegrita"/>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:u w:val="single"/>
        </w:rPr>
        <w:t>para</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rStyle w:val="Textoennegrita"/>
          <w:rFonts w:cstheme="minorHAnsi"/>
        </w:rPr>
        <w:t>:</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
        </w:rPr>
        <w:t xml:space="preserve"> mantener actualizada la base de datos y facilitar futuras compras y contactos comerciales.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Ttulo3"/>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
          <w:color w:val="4472C4" w:themeColor="accent1"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="es-GT"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
          <w:bCs w:val="0"/>
          <w:noProof/>
          <w:color w:val="4472C4" w:themeColor="accent1"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="es-GT"/>
        </w:rPr>
        <w:t>HUJR-002</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
          <w:b w:val="0"/>
          <w:bCs w:val="0"/>
          <w:noProof/>
          <w:color w:val="4472C4" w:themeColor="accent1"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="es-GT"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
          <w:color w:val="4472C4" w:themeColor="accent1"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="es-GT"/>
        </w:rPr>
        <w:t>Consulta de proveedores</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NormalWeb"/>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rStyle w:val="Textoennegrita"/>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
        </w:rPr>
        <w:t>Como</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
        </w:rPr>
        <w:t xml:space="preserve">: </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
        </w:rPr>
        <w:t>usuario autorizado,</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
        </w:rPr>
        <w:br/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rStyle w:val="Textoennegrita"/>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
        </w:rPr>
        <w:t>quiero</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
        </w:rPr>
        <w:t xml:space="preserve"> acceder a un listado de todos los proveedores registrados,</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
        </w:rPr>
        <w:br/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rStyle w:val="Textoennegrita"/>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
        </w:rPr>
        <w:t>para</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rStyle w:val="Textoennegrita"/>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
        </w:rPr>
        <w:t>:</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
        </w:rPr>
        <w:t xml:space="preserve"> poder visualizar su información de contacto y gestionar relaciones comerciales de manera eficiente.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Ttulo3"/>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
          <w:color w:val="4472C4" w:themeColor="accent1"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="es-GT"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
          <w:bCs w:val="0"/>
          <w:noProof/>
          <w:color w:val="4472C4" w:themeColor="accent1"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="es-GT"/>
        </w:rPr>
        <w:t>HUJR-003</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
          <w:color w:val="4472C4" w:themeColor="accent1"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="es-GT"/>
        </w:rPr>
        <w:t xml:space="preserve"> Actualización de datos de proveedor</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NormalWeb"/>
        <w:rPr>
          <w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rStyle w:val="Textoennegrita"/>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
        </w:rPr>
        <w:t>Como</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rStyle w:val="Textoennegrita"/>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
        </w:rPr>
        <w:t>:</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
        </w:rPr>
        <w:t xml:space="preserve"> administrador del sistema,</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
        </w:rPr>
        <w:br/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rStyle w:val="Textoennegrita"/>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
        </w:rPr>
        <w:t>quiero</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rStyle w:val="Textoennegrita"/>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
        </w:rPr>
        <w:t>:</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
        </w:rPr>
        <w:t xml:space="preserve"> editar la información de un proveedor previamente registrado,</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
        </w:rPr>
        <w:br/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rStyle w:val="Textoennegrita"/>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
        </w:rPr>
        <w:t>para</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rStyle w:val="Textoennegrita"/>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
        </w:rPr>
        <w:t>:</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
        </w:rPr>
        <w:t xml:space="preserve"> asegurarme de que los datos almacenados estén siempre actualizados y sean confiables</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:u w:val="single"/>
        </w:rPr>
      </w:pPr>
    </w:p>
'@

$r.InsertXML($xml)

Write-Host "Replaced paragraph at index $targetIndex; new paragraph count: $($d.Paragraphs.Count)"
